$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.488.95"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "1.910.58"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'239.22"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4780"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("D8").Value = "'0.2836"
$ws.Range("E8").Value = "  -3.48%  "
$ws.Range("D9").Value = "'0.06706"
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("D10").Value = "'18.76"
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("D11").Value = "'101.42"
$ws.Range("E11").Value = "  -4.50%  "
$ws.Range("D12").Value = "1.917.81"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "'0.07678"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "'5.202"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").Value = "'0.6693"
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("D16").Value = "30.509.36"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "'259.32"
$ws.Range("E17").Value = "  -6.22%  "
$ws.Range("D18").Value = "'0.9998"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'0.000007464"
$ws.Range("E19").Value = "  -3.36%  "
$ws.Range("D20").Value = "'12.65"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").Value = "'5.389"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'6.281"
$ws.Range("E23").Value = "  -3.95%  "
$ws.Range("D24").Value = "'9.353"
$ws.Range("E24").Value = "  -3.66%  "
$ws.Range("D25").Value = "'167.44"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "'19.13"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").Value = "'2.059"
$ws.Range("E27").Value = "  -4.76%  "
$ws.Range("D28").Value = "'4.817"
$ws.Range("E28").Value = "  +6.11%  "
$ws.Range("D29").Value = "'1.381"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "'0.1001"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("D31").Value = "'1.511"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").Value = "'4.255"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").Value = "'0.04720"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").Value = "'0.7255"
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("D35").Value = "'1.106"
$ws.Range("E35").Value = "  -4.35%  "
$ws.Range("D36").Value = "'2.711"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "'0.01916"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").Value = "'2.621"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'6.250"
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").Value = "'74.93"
$ws.Range("E40").Value = "  -4.23%  "
$ws.Range("D41").Value = "'1.969"
$ws.Range("E41").Value = "  -6.29%  "
$ws.Range("D42").Value = "'0.8614"
$ws.Range("E42").Value = "  -4.97%  "
$ws.Range("D43").Value = "'105.65"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "'0.4247"
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("D45").Value = "'0.9994"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "'7.380"
$ws.Range("E46").Value = "  -4.72%  "
$ws.Range("D47").Value = "'0.1201"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").Value = "'920.49"
$ws.Range("E48").Value = "  -6.72%  "
$ws.Range("D49").Value = "'34.73"
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("D50").Value = "'8.760"
$ws.Range("E50").Value = "  -5.24%  "
$ws.Range("E51").Value = "  +0.15%  "
